# Auto-generated update of Leve profit metrics (columns H-N) across all 8 sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), per scheduled price-refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 209.7
$ws.Range("I2").Value = 240
$ws.Range("J2").Value = 88.5
$ws.Range("K2").Value = 240
$ws.Range("L2").Value = 88.5
$ws.Range("M2").Value = -127
$ws.Range("N2").Value = -314.5
# Row 4
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
# Row 9
$ws.Range("H9").Value = 211.44444
$ws.Range("I9").Value = 200.57143
$ws.Range("K9").Value = 200.57143
$ws.Range("M9").Value = -31.57142999999999
# Row 17
$ws.Range("H17").Value = 3993.75
$ws.Range("J17").Value = 4825
$ws.Range("L17").Value = 14475
$ws.Range("N17").Value = -14811
# Row 28
$ws.Range("H28").Value = 1074.8077
$ws.Range("I28").Value = 1007.6
$ws.Range("K28").Value = 1007.6
$ws.Range("M28").Value = -522.6
# Row 104
$ws.Range("H104").Value = 180.5
$ws.Range("I104").Value = 180.5
$ws.Range("K104").Value = 541.5
$ws.Range("M104").Value = 1205.5
# Row 111
$ws.Range("H111").Value = 5168.1665
$ws.Range("I111").Value = 3155.5
$ws.Range("K111").Value = 9466.5
$ws.Range("M111").Value = -6399.5
# Row 112
$ws.Range("H112").Value = 1689.4546
$ws.Range("J112").Value = 1771.8
$ws.Range("L112").Value = 5315.4
$ws.Range("N112").Value = -7531.4
# Row 113
$ws.Range("H113").Value = 8682.272000000001
$ws.Range("I113").Value = 7151.25
$ws.Range("K113").Value = 7151.25
$ws.Range("M113").Value = -3897.25
# Row 116
$ws.Range("H116").Value = 4000
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 4000
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 4000
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -10884
# Row 127
$ws.Range("H127").Value = 10178.81
$ws.Range("I127").Value = 1184
$ws.Range("K127").Value = 3552
$ws.Range("M127").Value = 1408
# Row 137
$ws.Range("H137").Value = 6931.2354
$ws.Range("I137").Value = 10372.333
$ws.Range("J137").Value = 3060
$ws.Range("K137").Value = 31116.999
$ws.Range("L137").Value = 9180
$ws.Range("M137").Value = -28566.999
$ws.Range("N137").Value = -14280
# Row 141
$ws.Range("H141").Value = 5230.2144
$ws.Range("I141").Value = 3709.4614
$ws.Range("K141").Value = 11128.3842
$ws.Range("M141").Value = -5948.3842
$ws = $wb.Worksheets.Item("ARM")
# Row 30
$ws.Range("H30").Value = 894.3333
$ws.Range("I30").Value = 916.5
$ws.Range("J30").Value = 850
$ws.Range("K30").Value = 916.5
$ws.Range("L30").Value = 850
$ws.Range("M30").Value = -766.5
$ws.Range("N30").Value = -1150
# Row 32
$ws.Range("H32").Value = 43488630
$ws.Range("I32").Value = 43488630
$ws.Range("K32").Value = 43488630
$ws.Range("M32").Value = -43488343
# Row 74
$ws.Range("H74").Value = 2257.1875
$ws.Range("I74").Value = 2355.3845
$ws.Range("K74").Value = 2355.3845
$ws.Range("M74").Value = -1481.3845
# Row 77
$ws.Range("H77").Value = 2257.1875
$ws.Range("I77").Value = 2355.3845
$ws.Range("K77").Value = 11776.9225
$ws.Range("M77").Value = -7408.922500000001
# Row 106
$ws.Range("H106").Value = 30000
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
# Row 110
$ws.Range("H110").Value = 816
$ws.Range("I110").Value = 806.6667
$ws.Range("K110").Value = 806.6667
$ws.Range("M110").Value = 1238.3333
# Row 122
$ws.Range("H122").Value = 3374.0715
$ws.Range("I122").Value = 1884.3
$ws.Range("J122").Value = 7098.5
$ws.Range("K122").Value = 5652.9
$ws.Range("L122").Value = 21295.5
$ws.Range("M122").Value = -3202.9
$ws.Range("N122").Value = -26195.5
# Row 132
$ws.Range("H132").Value = 2529.641
$ws.Range("I132").Value = 1514.56
$ws.Range("J132").Value = 4342.2856
$ws.Range("K132").Value = 4543.68
$ws.Range("L132").Value = 13026.8568
$ws.Range("M132").Value = -2013.68
$ws.Range("N132").Value = -18086.8568
$ws = $wb.Worksheets.Item("BSM")
# Row 14
$ws.Range("H14").Value = 354.16666
$ws.Range("I14").Value = 400
$ws.Range("K14").Value = 400
$ws.Range("M14").Value = -228
# Row 16
$ws.Range("H16").Value = 8749.5
$ws.Range("J16").Value = 8749.5
$ws.Range("L16").Value = 8749.5
$ws.Range("N16").Value = -9089.5
# Row 54
$ws.Range("H54").Value = 30361
$ws.Range("I54").Value = 20541.5
$ws.Range("J54").Value = 50000
$ws.Range("K54").Value = 20541.5
$ws.Range("L54").Value = 50000
$ws.Range("M54").Value = -20057.5
$ws.Range("N54").Value = -50968
# Row 94
$ws.Range("H94").Value = 5954128.5
$ws.Range("I94").Value = 1368.3158
$ws.Range("J94").Value = 18521066
$ws.Range("K94").Value = 1368.3158
$ws.Range("L94").Value = 18521066
$ws.Range("M94").Value = -917.3158000000001
$ws.Range("N94").Value = -18521968
# Row 105
$ws.Range("H105").Value = 4098.0527
$ws.Range("I105").Value = 3547.077
$ws.Range("K105").Value = 3547.077
$ws.Range("M105").Value = -1800.077
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3238.1428
$ws.Range("I31").Value = 2226
$ws.Range("K31").Value = 2226
$ws.Range("M31").Value = -1931
# Row 34
$ws.Range("H34").Value = 3238.1428
$ws.Range("I34").Value = 2226
$ws.Range("K34").Value = 2226
$ws.Range("M34").Value = -2024
# Row 105
$ws.Range("H105").Value = 7753.769
$ws.Range("I105").Value = 1225.125
$ws.Range("K105").Value = 1225.125
$ws.Range("M105").Value = 521.875
# Row 122
$ws.Range("H122").Value = 368357.72
$ws.Range("I122").Value = 852118.5
$ws.Range("J122").Value = 5537.125
$ws.Range("K122").Value = 2556355.5
$ws.Range("L122").Value = 16611.375
$ws.Range("M122").Value = -2553905.5
$ws.Range("N122").Value = -21511.375
# Row 134
$ws.Range("H134").Value = 5625.5264
$ws.Range("I134").Value = 4158.316
$ws.Range("K134").Value = 12474.948
$ws.Range("M134").Value = -9939.948
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 321.42856
$ws.Range("I2").Value = 150
$ws.Range("K2").Value = 900
$ws.Range("M2").Value = -787
# Row 56
$ws.Range("H56").Value = 6019.727
$ws.Range("I56").Value = 6019.727
$ws.Range("K56").Value = 6019.727
$ws.Range("M56").Value = -5489.727
$ws = $wb.Worksheets.Item("GSM")
# Row 35
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
# Row 44
$ws.Range("H44").Value = 18731.25
$ws.Range("J44").Value = 18731.25
$ws.Range("L44").Value = 18731.25
$ws.Range("N44").Value = -19923.25
# Row 102
$ws.Range("H102").Value = 47488.23
$ws.Range("I102").Value = 53473.695
$ws.Range("K102").Value = 53473.695
$ws.Range("M102").Value = -51851.695
# Row 113
$ws.Range("H113").Value = 7652.9414
$ws.Range("I113").Value = 5268.6665
$ws.Range("J113").Value = 8953.454
$ws.Range("K113").Value = 5268.6665
$ws.Range("L113").Value = 8953.454
$ws.Range("M113").Value = -3098.6665
$ws.Range("N113").Value = -13293.454
$ws = $wb.Worksheets.Item("LTW")
# Row 33
$ws.Range("H33").Value = 25000
$ws.Range("I33").Value = 25000
$ws.Range("K33").Value = 25000
$ws.Range("M33").Value = -24710
# Row 44
$ws.Range("H44").Value = 23989
$ws.Range("J44").Value = 23989
$ws.Range("L44").Value = 23989
$ws.Range("N44").Value = -24901
# Row 55
$ws.Range("H55").Value = 2028.1666
$ws.Range("I55").Value = 4290
$ws.Range("J55").Value = 412.57144
$ws.Range("K55").Value = 4290
$ws.Range("L55").Value = 412.57144
$ws.Range("M55").Value = -4117
$ws.Range("N55").Value = -758.5714399999999
# Row 68
$ws.Range("H68").Value = 211616.25
$ws.Range("I68").Value = 135733.14
$ws.Range("J68").Value = 338088.12
$ws.Range("K68").Value = 135733.14
$ws.Range("L68").Value = 338088.12
$ws.Range("M68").Value = -134984.14
$ws.Range("N68").Value = -339586.12
# Row 71
$ws.Range("H71").Value = 211616.25
$ws.Range("I71").Value = 135733.14
$ws.Range("J71").Value = 338088.12
$ws.Range("K71").Value = 678665.7000000001
$ws.Range("L71").Value = 1690440.6
$ws.Range("M71").Value = -674921.7000000001
$ws.Range("N71").Value = -1697928.6
# Row 100
$ws.Range("H100").Value = 3972.0476
$ws.Range("J100").Value = 5392
$ws.Range("L100").Value = 5392
$ws.Range("N100").Value = -6474
# Row 136
$ws.Range("H136").Value = 6731.615
$ws.Range("I136").Value = 4452.4
$ws.Range("J136").Value = 8156.125
$ws.Range("K136").Value = 13357.2
$ws.Range("L136").Value = 24468.375
$ws.Range("M136").Value = -10807.2
$ws.Range("N136").Value = -29568.375
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 15875995
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 18521494
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 18521494
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -18522742
# Row 65
$ws.Range("H65").Value = 15875995
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 18521494
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 92607470
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -92613710
# Row 122
$ws.Range("H122").Value = 259811.89
$ws.Range("I122").Value = 372849.47
$ws.Range("K122").Value = 1118548.41
$ws.Range("M122").Value = -1116098.41
# Row 126
$ws.Range("H126").Value = 2391.8
$ws.Range("I126").Value = 2391.8
$ws.Range("K126").Value = 7175.400000000001
$ws.Range("M126").Value = -4705.400000000001
# Row 132
$ws.Range("H132").Value = 2020.258
$ws.Range("I132").Value = 1028.9048
$ws.Range("K132").Value = 3086.7144
$ws.Range("M132").Value = -556.7143999999998
